# Added inserting of values from scenario to calc
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# Qualify the sheet reference of every defined name with single quotes,
# e.g. general!$C$2  ->  'general'!$C$2
$sheetName = $ws.Name
$unquoted = $sheetName + "!"
$quoted = "'" + $sheetName + "'!"
foreach ($n in $wb.Names) {
    $formula = $n.RefersTo
    if ($formula.Contains($unquoted) -and -not $formula.Contains($quoted)) {
        $n.RefersTo = $formula.Replace($unquoted, $quoted)
    }
}

# Insert the updated values for the "general_row_cells" named range
# (row 11, columns C:G) coming from the scenario.
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = 10

# Touch the header/footer so it is (re)written to the sheet, matching the
# calc-load refresh that happens when the scenario values are inserted.
$ws.PageSetup.CenterHeader = $ws.PageSetup.CenterHeader
